$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,12

$data[0,0] = 6.962736285065886
$data[0,1] = 8.333005155668802
$data[0,2] = 10.89460059637521
$data[0,3] = 41.54024229733263
$data[0,4] = 3.730130508428733
$data[0,5] = 0
$data[0,6] = 35.75094196978976
$data[0,7] = 0
$data[0,8] = 20.63716289012601
$data[0,9] = 9.253837271770223
$data[0,10] = 20.45957516007954
$data[0,11] = 20.60365822414531
$data[1,0] = 6.957063275051402
$data[1,1] = 8.364677499272123
$data[1,2] = 10.91201797578635
$data[1,3] = 41.44953793057156
$data[1,4] = 3.733962820756238
$data[1,5] = 0
$data[1,6] = 35.73973347207001
$data[1,7] = 0
$data[1,8] = 20.32930343580639
$data[1,9] = 9.275345362419166
$data[1,10] = 20.35025438267667
$data[1,11] = 20.68047000147437
$data[2,0] = 6.954503110190878
$data[2,1] = 8.384967626201686
$data[2,2] = 10.92419788562953
$data[2,3] = 41.40564890347515
$data[2,4] = 3.736437294658621
$data[2,5] = 0
$data[2,6] = 35.74196866166558
$data[2,7] = 0
$data[2,8] = 20.14350878323259
$data[2,9] = 9.289512923223898
$data[2,10] = 20.28749680684668
$data[2,11] = 20.72966090189789
$data[3,0] = 6.953693045897992
$data[3,1] = 8.393448817952228
$data[3,2] = 10.92953485701097
$data[3,3] = 41.39073568354879
$data[3,4] = 3.737476311848971
$data[3,5] = 0
$data[3,6] = 35.74516762702646
$data[3,7] = 0
$data[3,8] = 20.06870523736269
$data[3,9] = 9.295528415392575
$data[3,10] = 20.26303980639583
$data[3,11] = 20.75021833382243
$data[4,0] = 6.953572651944779
$data[4,1] = 8.394869988775264
$data[4,2] = 10.93044361855896
$data[4,3] = 41.38843889693027
$data[4,4] = 3.737650694428408
$data[4,5] = 0
$data[4,6] = 35.74583684768129
$data[4,7] = 0
$data[4,8] = 20.0563418743265
$data[4,9] = 9.296541914696466
$data[4,10] = 20.2590467518976
$data[4,11] = 20.75366283877761
$data[5,0] = 6.954491239661125
$data[5,1] = 8.385081143906007
$data[5,2] = 10.92426834938479
$data[5,3] = 41.40543574309999
$data[5,4] = 3.736451182957266
$data[5,5] = 0
$data[5,6] = 35.74200254682084
$data[5,7] = 0
$data[5,8] = 20.142496148483
$data[5,9] = 9.289593069599986
$data[5,10] = 20.28716242326237
$data[5,11] = 20.72993607235497
$data[6,0] = 6.960589334345663
$data[6,1] = 8.343751246966677
$data[6,2] = 10.90029778705306
$data[6,3] = 41.50652075692855
$data[6,4] = 3.731426760030707
$data[6,5] = 0
$data[6,6] = 35.74518272782368
$data[6,7] = 0
$data[6,8] = 20.53040037445565
$data[6,9] = 9.261053891223584
$data[6,10] = 20.42098735949875
$data[6,11] = 20.62972303664133
$data[7,0] = 6.97982033080973
$data[7,1] = 8.269357701692959
$data[7,2] = 10.86507861437645
$data[7,3] = 41.79809932806823
$data[7,4] = 3.722531933111358
$data[7,5] = 0
$data[7,6] = 35.82389875715778
$data[7,7] = 0
$data[7,8] = 21.31243877325786
$data[7,9] = 9.212703669309018
$data[7,10] = 20.71713238427851
$data[7,11] = 20.44921674725267
$data[8,0] = 6.998309298395718
$data[8,1] = 8.218706148315635
$data[8,2] = 10.84638796417304
$data[8,3] = 42.06861540937015
$data[8,4] = 3.716573465837541
$data[8,5] = 0
$data[8,6] = 35.92601459287949
$data[8,7] = 0
$data[8,8] = 21.89428003488826
$data[8,9] = 9.181803913142366
$data[8,10] = 20.95392410615182
$data[8,11] = 20.32624317830462
$data[9,0] = 7.007649517288478
$data[9,1] = 8.19652221715484
$data[9,2] = 10.83944448108871
$data[9,3] = 42.20370514779768
$data[9,4] = 3.713986398502107
$data[9,5] = 0
$data[9,6] = 35.9820730403685
$data[9,7] = 0
$data[9,8] = 22.15935729822012
$data[9,9] = 9.168746723946992
$data[9,10] = 21.06550600536441
$data[9,11] = 20.27236853205498
$data[10,0] = 7.011318401043123
$data[10,1] = 8.188244261385577
$data[10,2] = 10.83703920701365
$data[10,3] = 42.25656816425465
$data[10,4] = 3.713024375028465
$data[10,5] = 0
$data[10,6] = 36.0046781782039
$data[10,7] = 0
$data[10,8] = 22.25969481509762
$data[10,9] = 9.163945708082849
$data[10,10] = 21.10828671855894
$data[10,11] = 20.2522629353593
$data[11,0] = 7.01052240128698
$data[11,1] = 8.190021625232893
$data[11,2] = 10.83754726228036
$data[11,3] = 42.24510761881842
$data[11,4] = 3.713230781106789
$data[11,5] = 0
$data[11,6] = 35.99974859903617
$data[11,7] = 0
$data[11,8] = 22.23808871276523
$data[11,9] = 9.164973316079431
$data[11,10] = 21.0990501500786
$data[11,11] = 20.25657991576139
$data[12,0] = 7.007948724336092
$data[12,1] = 8.195838731425548
$data[12,2] = 10.83924210775691
$data[12,3] = 42.20802013969855
$data[12,4] = 3.713906899304641
$data[12,5] = 0
$data[12,6] = 35.9839052112986
$data[12,7] = 0
$data[12,8] = 22.1676135096268
$data[12,9] = 9.168348867689017
$data[12,10] = 21.0690152042164
$data[12,11] = 20.27070851823889
$data[13,0] = 7.006389403584691
$data[13,1] = 8.19941782405076
$data[13,2] = 10.84030942650806
$data[13,3] = 42.1855246295321
$data[13,4] = 3.714323335517668
$data[13,5] = 0
$data[13,6] = 35.97437985345945
$data[13,7] = 0
$data[13,8] = 22.12443714500295
$data[13,9] = 9.170435166152673
$data[13,10] = 21.05068570803654
$data[13,11] = 20.27940113831325
$data[14,0] = 6.997717452552862
$data[14,1] = 8.220173118050596
$data[14,2] = 10.8468730995583
$data[14,3] = 42.06002720874146
$data[14,4] = 3.716745011711378
$data[14,5] = 0
$data[14,6] = 35.92254400409373
$data[14,7] = 0
$data[14,8] = 21.87695608950102
$data[14,9] = 9.182677318728054
$data[14,10] = 20.9467073012599
$data[14,11] = 20.32980546074249
$data[15,0] = 6.992634367937507
$data[15,1] = 8.233124985226883
$data[15,2] = 10.8512989295541
$data[15,3] = 41.9861048891383
$data[15,4] = 3.718262175274033
$data[15,5] = 0
$data[15,6] = 35.89320209473446
$data[15,7] = 0
$data[15,8] = 21.72516476178076
$data[15,9] = 9.190443249632214
$data[15,10] = 20.88388958870723
$data[15,11] = 20.3612550426637
$data[16,0] = 6.989798309315823
$data[16,1] = 8.240655331584081
$data[16,2] = 10.85399129809954
$data[16,3] = 41.94472065232803
$data[16,4] = 3.719146435813644
$data[16,5] = 0
$data[16,6] = 35.87722987427549
$data[16,7] = 0
$data[16,8] = 21.6379012722197
$data[16,9] = 9.195004076928011
$data[16,10] = 20.84812399224832
$data[16,11] = 20.37953862916246
$data[17,0] = 6.988853164172652
$data[17,1] = 8.243218872674991
$data[17,2] = 10.85492809463058
$data[17,3] = 41.9309040494288
$data[17,4] = 3.719447831807551
$data[17,5] = 0
$data[17,6] = 35.87197739977505
$data[17,7] = 0
$data[17,8] = 21.60836575800984
$data[17,9] = 9.196564458403067
$data[17,10] = 20.83607800431925
$data[17,11] = 20.38576261280756
$data[18,0] = 6.993166417111284
$data[18,1] = 8.231737882201132
$data[18,2] = 10.85081260518094
$data[18,3] = 41.99385684264372
$data[18,4] = 3.71809946796217
$data[18,5] = 0
$data[18,6] = 35.89623200725283
$data[18,7] = 0
$data[18,8] = 21.74131949230777
$data[18,9] = 9.189606818942607
$data[18,10] = 20.8905390186651
$data[18,11] = 20.35788705044196
$data[19,0] = 7.008701108817084
$data[19,1] = 8.194126784238236
$data[19,2] = 10.83873821035298
$data[19,3] = 42.21886748582779
$data[19,4] = 3.713707829070884
$data[19,5] = 0
$data[19,6] = 35.9885214683711
$data[19,7] = 0
$data[19,8] = 22.18831567912029
$data[19,9] = 9.167353494781755
$data[19,10] = 21.07782312384259
$data[19,11] = 20.26655059656878
$data[20,0] = 7.019622145536583
$data[20,1] = 8.170260135664417
$data[20,2] = 10.83215288257432
$data[20,3] = 42.37586478848563
$data[20,4] = 3.710940423133108
$data[20,5] = 0
$data[20,6] = 36.05686115731304
$data[20,7] = 0
$data[20,8] = 22.48016836896066
$data[20,9] = 9.153645740943514
$data[20,10] = 21.20328238988622
$data[20,11] = 20.20857896526065
$data[21,0] = 7.013723681689575
$data[21,1] = 8.182933086185397
$data[21,2] = 10.83554814370934
$data[21,3] = 42.29117109860072
$data[21,4] = 3.712408071909064
$data[21,5] = 0
$data[21,6] = 36.01965468467814
$data[21,7] = 0
$data[21,8] = 22.32445883877879
$data[21,9] = 9.160885401836149
$data[21,10] = 21.13605226531048
$data[21,11] = 20.23936248289766
$data[22,0] = 6.992925608825293
$data[22,1] = 8.232364729330945
$data[22,2] = 10.8510320117467
$data[22,3] = 41.99034871078391
$data[22,4] = 3.718172990479705
$data[22,5] = 0
$data[22,6] = 35.89485939006438
$data[22,7] = 0
$data[22,8] = 21.73401592066361
$data[22,9] = 9.189984669939973
$data[22,10] = 20.88753172172124
$data[22,11] = 20.35940908791704
$data[23,0] = 6.973846459853822
$data[23,1] = 8.288776057794129
$data[23,2] = 10.87334439962501
$data[23,3] = 41.7092763457697
$data[23,4] = 3.724836429557966
$data[23,5] = 0
$data[23,6] = 35.79482803501289
$data[23,7] = 0
$data[23,8] = 21.09918589062336
$data[23,9] = 9.224970447881827
$data[23,10] = 20.63354065474708
$data[23,11] = 20.49634644823529

$ws.Range("C2:N25").Value2 = $data

Write-Output "Updated loading_percent values for C2:N25"
